$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M18").Value = "plantarflexors"
$ws.Range("M19").Value = "Glut. Med."
$ws.Range("M20").Value = "Hip Flexors"
$ws.Range("M21").Value = "Hip Adductors"
$ws.Range("M22").Value = "Glut. Min."
$ws.Range("M23").Value = "Glut. Max."
$ws.Range("M24").Value = "Hamstrings"
$ws.Range("M25").Value = "Quads"
$ws.Range("M26").Value = "Dorsiflexors"

$ws.Range("M27").Select()
